# #5: property boat&car done
#
# The "汽車" (car) property sheet was missing a column describing the
# engine capacity (c.c.) and the common trailing metadata columns
# (property_category, category, date, legislator_name, legislator_id,
# source_file, index) that all the other property sheets already have.
# Row 1 had also mistakenly been filled with the first data row's values
# instead of column headers - fix that at the same time.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("汽車")

# Broadcast the header-row formatting (bold + border, style of B1) onto the
# new header cells H1:N1 before we fill in their text.
$ws.Range("B1").Copy($ws.Range("H1:N1"))

# --- Row 1: turn into real column headers ---------------------------------
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "capacity"
$ws.Range("D1").Value = "owner"
$ws.Range("E1").Value = "register_date"
$ws.Range("F1").Value = "register_reason"
$ws.Range("G1").Value = "acquire_value"
$ws.Range("H1").Value = "property_category"
$ws.Range("I1").Value = "category"
$ws.Range("J1").Value = "date"
$ws.Range("K1").Value = "legislator_name"
$ws.Range("L1").Value = "legislator_id"
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

# --- Row 2: car record #44 --------------------------------------------------
$ws.Range("A2").Value = 44
$ws.Range("B2").Value = "國瑞NV1EPE"
$ws.Range("C2").Value = 1998
$ws.Range("D2").Value = "趙哲宏"
$ws.Range("E2").Value = "94年04月12日"
$ws.Range("F2").Value = "買賣"
$ws.Range("G2").Value = "(超過五年）"
$ws.Range("H2").Value = "land"
$ws.Range("I2").Value = "normal"
# "date" must stay a plain text value, not get auto-parsed into a date
# serial number.
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "2011-11-22"
$ws.Range("K2").Value = "葉宜津"
$ws.Range("L2").Value = 855
$ws.Range("M2").Value = "tmp14431"
$ws.Range("N2").Value = 44

# --- Row 3: car record #45 --------------------------------------------------
$ws.Range("A3").Value = 45
$ws.Range("B3").Value = "國瑞NV1EPE"
$ws.Range("C3").Value = 1998
$ws.Range("D3").Value = "趙哲宏"
$ws.Range("E3").Value = "93年04月15日"
$ws.Range("F3").Value = "買賣"
$ws.Range("G3").Value = "(超過五年）"
$ws.Range("H3").Value = "land"
$ws.Range("I3").Value = "normal"
$ws.Range("J3").NumberFormat = "@"
$ws.Range("J3").Value = "2011-11-22"
$ws.Range("K3").Value = "葉宜津"
$ws.Range("L3").Value = 855
$ws.Range("M3").Value = "tmp14431"
$ws.Range("N3").Value = 45
